$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (iteration 0) - update xn/fxn/E, keep Iteration (A2) as-is
$ws.Range("B2").Value = "'0.7"
$ws.Range("C2").Value = "'0.543752707470477"
$ws.Range("D2").Value = "'1.000001"

# Row 3 (iteration 1) - keep Iteration (A3) as-is, fill in the real xn/fxn/E
$ws.Range("B3").Value = "'0.948715074082467"
$ws.Range("C3").Value = "'-0.117791527780404"
$ws.Range("D3").Value = "'0.262159926491109"

# Row 4 (iteration 2) - new row
$ws.Range("A4").Value = "'2"
$ws.Range("B4").Value = "'0.910838782233832"
$ws.Range("C4").Value = "'-0.002474648500744"
$ws.Range("D4").Value = "'0.0415839691803017"

# Row 5 (iteration 3) - new row
$ws.Range("A5").Value = "'3"
$ws.Range("B5").Value = "'0.910007980114702"
$ws.Range("C5").Value = "'-1.21297462740699e-06"
$ws.Range("D5").Value = "'0.0009129613555964"

# Row 6 (iteration 4) - new row
$ws.Range("A6").Value = "'4"
$ws.Range("B6").Value = "'0.910007572488831"
$ws.Range("C6").Value = "'-3.62376795237651e-13"
$ws.Range("D6").Value = "'4.47936790664828e-07"
